# Applies the repositioning edits captured in the diff:
#   Slide 4, Shape 1 ("標題 1" / 遇到的困難與需求): off.x 343224 -> 501674 EMU
#   Slide 4, Shape 2 ("標題 1" / 如何解決與改善)  : off.x 6798330 -> 6709840 EMU
#   Slide 9, Shape 2 ("文字方塊 1" text box)      : off.x 2412087 -> 3277327 EMU
#                                                   off.y 2256702 -> 2441824 EMU
#                                                   ext.cx 7317129 -> 5935500 EMU
#
# PowerPoint's COM object model expresses Shape.Left/Top/Width/Height in
# points (a [Single]/f32 under the hood), while the OOXML stores EMUs
# (914400 EMU per inch, 12700 EMU per point). A plain emu/12700 division
# can land a hair below the intended EMU once round-tripped through the
# f32 Left/Top/Width/Height properties (the host truncates when it turns
# the point value back into EMU), so nudge by half an EMU worth of points
# before converting to land exactly on the intended integer EMU value.

function EMUToPoints($emu) {
    return ($emu + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation

# --- Slide 4 ---
$s4 = $p.Slides.Item(4)

$shape1 = $s4.Shapes.Item(1)
$shape1.Left = EMUToPoints 501674

$shape2 = $s4.Shapes.Item(2)
$shape2.Left = EMUToPoints 6709840

# --- Slide 9 ---
$s9 = $p.Slides.Item(9)

$shape = $s9.Shapes.Item(2)
$shape.Left = EMUToPoints 3277327
$shape.Top = EMUToPoints 2441824
$shape.Width = EMUToPoints 5935500
